# Apply cell value updates to match the target OOXML diff.
# Cells whose new text value "looks like a number" (e.g. "12.50") would be
# auto-coerced into a numeric cell by the COM Value setter, losing the exact
# text formatting (trailing zeros, float precision). For those we temporarily
# force a text number format, assign the value, then restore the default style
# so the cell keeps looking/serializing exactly like its untouched neighbors.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.561.44'
$ws.Range('E2').Value = '  +3.44%  '
$ws.Range('D3').Value = '2.433.17'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.514'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.99%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('E14').Value = '  +2.54%  '
$ws.Range('D15').Value = '2.813.66'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('D16').Value = '2.430.92'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.837'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.12%  '
$ws.Range('D18').Value = '44.429.57'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.78%  '
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').Value = '0.0₃0910'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('E28').Value = '  -2.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.122'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +15.28%  '
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0764'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.44%  '
$ws.Range('E36').Value = '  +2.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.50'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '126.12'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0290'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '1.949.64'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('E46').Value = '  +7.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.87%  '
$ws.Range('E48').Value = '  +10.05%  '
$ws.Range('D49').Value = '2.676.48'
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.04%  '
